# Generate Report for Archive
# Update the localization status from "Ready for handoff" to "In Translation"
# on every worksheet where it appears, then autofit the affected columns so
# their widths shrink to match the new (shorter) text.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws1.Columns.Item(5).EntireColumn.AutoFit()
$ws1.Columns.Item(6).EntireColumn.AutoFit()

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "In Translation"
$ws2.Columns.Item(3).EntireColumn.AutoFit()

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "In Translation"
$ws3.Columns.Item(3).EntireColumn.AutoFit()
